$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fill in the next blank row of the revision-history table with the
#    new "02/13/2015" entry (date / description / author).
# ------------------------------------------------------------------
$found = $false
foreach ($tbl in $d.Tables) {
    if ($found) { break }
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        if ($found) { break }
        $row = $tbl.Rows.Item($r)
        if ($row.Cells.Count -ge 3) {
            $c1 = $row.Cells.Item(1).Range.Text
            $c2 = $row.Cells.Item(2).Range.Text
            $c3 = $row.Cells.Item(3).Range.Text
            # Strip trailing cell-mark / paragraph-mark control chars for
            # an "is this cell empty" test.
            $t1 = $c1 -replace "[\x07\x0d]", ""
            $t2 = $c2 -replace "[\x07\x0d]", ""
            $t3 = $c3 -replace "[\x07\x0d]", ""
            if ($t1 -eq "" -and $t2 -eq "" -and $t3 -eq "" -and $r -gt 1) {
                $found = $true

                # --- Date cell -------------------------------------------------
                $rng1 = $row.Cells.Item(1).Range
                $rng1.Collapse(0)
                $rng1.MoveEnd(1, -1)
                $rng1.Text = "02/13/2015"

                # --- Description cell ------------------------------------------
                $rng2 = $row.Cells.Item(2).Range
                $rng2.Collapse(0)
                $rng2.MoveEnd(1, -1)
                $rng2.Text = "P13276 – eCL Change BCC to CCO`rChanged Vangent to GDIT in footer"

                # --- Author cell -------------------------------------------------
                $rng3 = $row.Cells.Item(3).Range
                $rng3.Collapse(0)
                $rng3.MoveEnd(1, -1)
                $rng3.Text = "Doug Stearns"
            }
        }
    }
}

# ------------------------------------------------------------------
# 2. Update the footer text:
#    "VANGENT PROPRIETARY - CONFIDENTIAL" -> "GDIT, INC. CONFIDENTIAL"
#    "Copyright © 2011 Vangent, " -> "Copyright © 2011, "
# ------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $ftr = $sec.Footers.Item($idx)
        if ($ftr.Exists) {
            $rng = $ftr.Range
            $rng.Find.Execute("VANGENT PROPRIETARY - CONFIDENTIAL", $true, $false, $false, $false, $false, $true, 1, $false, "GDIT, INC. CONFIDENTIAL", 2)
            $rng2 = $ftr.Range
            $rng2.Find.Execute("Copyright © 2011 Vangent, ", $true, $false, $false, $false, $false, $true, 1, $false, "Copyright © 2011, ", 2)
        }
    }
}
